# Generate Report for Handoff
# The 73154f29-4e97-4bc8-84b8-7c98aca212a3.md file has finished translation
# and is ready to be handed off, so its status moves from "In Translation"
# to "Ready for handoff" (with an updated handoff timestamp/priority) on
# every sheet that tracks it: the Overview roll-up plus the per-locale
# (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 73154f29...md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-20 12:11:16"

# --- zh-cn detail sheet: row for 73154f29...md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-20 12:11:13"

# --- de-de detail sheet: row for 73154f29...md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-20 12:11:16"

# The longer "Ready for handoff" status text widens the Status-like
# columns on each sheet (mirrors the author's column autofit after the
# longer text was written).
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
